# Change the year in the astromap link from 2018 to 2022.
# The original sentence was split across three differently-formatted runs
# (plain text, a Hyperlink-styled URL run, and a trailing plain run).
# Locate the whole sentence, remove it and retype it as a single run with
# the updated year and no inherited character formatting (i.e. the
# Hyperlink styling on the URL is dropped, matching the merged plain-text
# run produced by the edit).

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()

$oldText = "Mapky v tomto dokumente pripravil Jan Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/)."
$newText = "Mapky v tomto dokumente pripravil Jan Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

$found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)

if ($found) {
    $rng = $find.Parent
    $rng.Delete()
    $rng.InsertAfter($newText)
}
